$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.94933966666666
$ws.Range("H2").Value = 71.848019
$ws.Range("I2").Value = 0.003771463579284448
$ws.Range("J2").Value = 0.003771463579284448
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8180823333333334
$ws.Range("N2").Value = 2.454247
$ws.Range("O2").Value = 0.5115352725808422
$ws.Range("P2").Value = 0.5115352725808422
$ws.Range("Q2").Value = 19.59253167629922
$ws.Range("R2").Value = 176.332785086693
$ws.Range("S2").Value = 0.001929236650057989
$ws.Range("T2").Value = 0.001929236650057989

$ws.Range("G3").Value = 23.94933966666666
$ws.Range("H3").Value = 71.848019
$ws.Range("I3").Value = 0.003771463579284448
$ws.Range("J3").Value = 0.003771463579284448
$ws.Range("M3").Value = 0.7811863333333333
$ws.Range("N3").Value = 2.343559
$ws.Range("O3").Value = 0.4884647274191578
$ws.Range("P3").Value = 0.4884647274191579
$ws.Range("Q3").Value = 18.70889683995789
$ws.Range("R3").Value = 168.380071559621
$ws.Range("S3").Value = 0.001842226929226459
$ws.Range("T3").Value = 0.001842226929226459

$ws.Range("I4").Value = 0.9508087878751565
$ws.Range("J4").Value = 0.9508087878751567
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8180823333333334
$ws.Range("N4").Value = 2.454247
$ws.Range("O4").Value = 0.5115352725808422
$ws.Range("P4").Value = 0.5115352725808422
$ws.Range("Q4").Value = 4939.395781751674
$ws.Range("R4").Value = 44454.56203576506
$ws.Range("S4").Value = 0.4863722324779784
$ws.Range("T4").Value = 0.4863722324779785

$ws.Range("I5").Value = 0.9508087878751565
$ws.Range("J5").Value = 0.9508087878751567
$ws.Range("M5").Value = 0.7811863333333333
$ws.Range("N5").Value = 2.343559
$ws.Range("O5").Value = 0.4884647274191578
$ws.Range("P5").Value = 0.4884647274191579
$ws.Range("Q5").Value = 4716.626093007823
$ws.Range("R5").Value = 42449.63483707041
$ws.Range("S5").Value = 0.4644365553971782
$ws.Range("T5").Value = 0.4644365553971783

$ws.Range("G6").Value = 285.3476563333333
$ws.Range("H6").Value = 856.042969
$ws.Range("I6").Value = 0.0449356144375536
$ws.Range("J6").Value = 0.04493561443755361
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8180823333333334
$ws.Range("N6").Value = 2.454247
$ws.Range("O6").Value = 0.5115352725808422
$ws.Range("P6").Value = 0.5115352725808422
$ws.Range("Q6").Value = 233.4378765043714
$ws.Range("R6").Value = 2100.940888539343
$ws.Range("S6").Value = 0.02298615177990161
$ws.Range("T6").Value = 0.02298615177990162

$ws.Range("G7").Value = 285.3476563333333
$ws.Range("H7").Value = 856.042969
$ws.Range("I7").Value = 0.0449356144375536
$ws.Range("J7").Value = 0.04493561443755361
$ws.Range("M7").Value = 0.7811863333333333
$ws.Range("N7").Value = 2.343559
$ws.Range("O7").Value = 0.4884647274191578
$ws.Range("P7").Value = 0.4884647274191579
$ws.Range("Q7").Value = 222.9096893762967
$ws.Range("R7").Value = 2006.187204386671
$ws.Range("S7").Value = 0.02194946265765199
$ws.Range("T7").Value = 0.021949462657652

$ws.Range("G8").Value = 3.074321666666667
$ws.Range("H8").Value = 9.222965
$ws.Range("I8").Value = 0.0004841341080053326
$ws.Range("J8").Value = 0.0004841341080053326
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8180823333333334
$ws.Range("N8").Value = 2.454247
$ws.Range("O8").Value = 0.5115352725808422
$ws.Range("P8").Value = 0.5115352725808422
$ws.Range("Q8").Value = 2.515048242483889
$ws.Range("R8").Value = 22.635434182355
$ws.Range("S8").Value = 0.0002476516729041907
$ws.Range("T8").Value = 0.0002476516729041907

$ws.Range("G9").Value = 3.074321666666667
$ws.Range("H9").Value = 9.222965
$ws.Range("I9").Value = 0.0004841341080053326
$ws.Range("J9").Value = 0.0004841341080053326
$ws.Range("M9").Value = 0.7811863333333333
$ws.Range("N9").Value = 2.343559
$ws.Range("O9").Value = 0.4884647274191578
$ws.Range("P9").Value = 0.4884647274191579
$ws.Range("Q9").Value = 2.401618070270556
$ws.Range("R9").Value = 21.614562632435
$ws.Range("S9").Value = 0.0002364824351011419
$ws.Range("T9").Value = 0.000236482435101142
